$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.141.31"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.802.49"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.83"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5095"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3848"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07732"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.101"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.68"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.349"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.001"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.37"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.800.37"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.290"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.21"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06565"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.27"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.968"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.152.45"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.07"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.244"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.62"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.431"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.006.28"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.28"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.59"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.55%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.048"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.651"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.550"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07024"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.026"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02348"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2170"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.027"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.49"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -6.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6134"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.17%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.153"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.21"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5916"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.295"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.713"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.49"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.199"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.917"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06735"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.36%  "
